# Commit: "updated Deck and Text to branch"
# The title placeholder ("Title 1", shape id=2 / ctrTitle) on slide 1 was
# empty and gets the new title text, typed as three runs (as PowerPoint
# itself splits a typed sentence containing a flagged/misspelled word
# into separate runs around the misspelling):
#   "This is the new " + "editSlides" + " Deck!"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Title

$tr = $titleShape.TextFrame.TextRange
$tr.Text = "This is the new "
[void]$tr.InsertAfter("editSlides")
[void]$tr.InsertAfter(" Deck!")
